# TRAD-43 fixed reading candles before first
#
# The FirstCandleDate column (B) previously stored only whole-day serial
# dates. The fix records the actual first-candle timestamp (date + time of
# day), so every NumberFormat in that column needs an "h:mm:ss" time
# portion added, and most rows get a non-zero time-of-day fraction added
# to their stored serial value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (new value, date format that row should use)
# Format families (same ones already used in the sheet, now with a time part):
#   "dd.mm.yyyy h:mm:ss"  (was "dd.mm.yyyy")
#   "d.m.yyyy h:mm:ss"    (was "d.m.yyyy")
#   "dd.mm.yy h:mm:ss"    (was "dd.mm.yy")

$rows = @(
    @{ Row = 2;  Value = 43839.336805555555; Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 3;  Value = 44125.39236111111;  Format = "d.m.yyyy h:mm:ss"   },
    @{ Row = 4;  Value = 43796.322916666664; Format = "d.m.yyyy h:mm:ss"   },
    @{ Row = 5;  Value = 44270.0;            Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 6;  Value = 43916.44097222222;  Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 7;  Value = 43716.74652777778;  Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 8;  Value = 44329.395833333336; Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 9;  Value = 43836.34722222222;  Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 10; Value = 44484.0;            Format = "d.m.yyyy h:mm:ss"   },
    @{ Row = 11; Value = 44088.291666666664; Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 12; Value = 44376.302083333336; Format = "dd.mm.yy h:mm:ss"   },
    @{ Row = 13; Value = 43871.333333333336; Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 14; Value = 44273.32638888889;  Format = "dd.mm.yyyy h:mm:ss" },
    @{ Row = 15; Value = 43861.333333333336; Format = "dd.mm.yyyy h:mm:ss" }
)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, 2)
    $cell.Value2 = $r.Value
    $cell.NumberFormat = $r.Format
}

# Column B ("FirstCandleDate") needs to be widened now that it shows a
# full date + time value instead of just a date.
$ws.Columns.Item(2).ColumnWidth = 20.13
